$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value2 = 0.04675869436393043
$ws.Range("D2").Value2 = 0.02109931985517299
$ws.Range("E2").Value2 = 0.1072129439688219
$ws.Range("F2").Value2 = 4.213970251737038
$ws.Range("G2").Value2 = 0.002597538318608457
$ws.Range("J2").Value2 = 0.2415610514328819
$ws.Range("K2").Value2 = 4.244402053868555
$ws.Range("N2").Value2 = 2.120912810895483

$ws.Range("C3").Value2 = 0.04162138080707223
$ws.Range("D3").Value2 = 0.02090468328288608
$ws.Range("E3").Value2 = 0.1054105674011048
$ws.Range("F3").Value2 = 4.164899039714925
$ws.Range("G3").Value2 = 0.002604782024412401
$ws.Range("J3").Value2 = 0.236741573115296
$ws.Range("K3").Value2 = 4.036320840420501
$ws.Range("N3").Value2 = 2.141224808104731

$ws.Range("C4").Value2 = 0.03848973142326884
$ws.Range("D4").Value2 = 0.02078927943432518
$ws.Range("E4").Value2 = 0.1043644044330776
$ws.Range("F4").Value2 = 4.137388990938092
$ws.Range("G4").Value2 = 0.002609456241278159
$ws.Range("J4").Value2 = 0.2339371624792008
$ws.Range("K4").Value2 = 3.911394907304725
$ws.Range("N4").Value2 = 2.15444067618818

$ws.Range("C5").Value2 = 0.0372190620329178
$ws.Range("D5").Value2 = 0.02074329842788281
$ws.Range("E5").Value2 = 0.1039532441808966
$ws.Range("F5").Value2 = 4.126832491864633
$ws.Range("G5").Value2 = 0.002611418217697417
$ws.Range("J5").Value2 = 0.2328329746568798
$ws.Range("K5").Value2 = 3.861192595878094
$ws.Range("N5").Value2 = 2.16001258286596

$ws.Range("C6").Value2 = 0.0370083949280513
$ws.Range("D6").Value2 = 0.02073572700250104
$ws.Range("E6").Value2 = 0.1038858853882267
$ws.Range("F6").Value2 = 4.12511896088418
$ws.Range("G6").Value2 = 0.002611747463670276
$ws.Range("J6").Value2 = 0.2326519501665913
$ws.Range("K6").Value2 = 3.852898971682976
$ws.Range("N6").Value2 = 2.160949024367746

$ws.Range("C7").Value2 = 0.03847257271254989
$ws.Range("D7").Value2 = 0.0207886550588583
$ws.Range("E7").Value2 = 0.1043587980715088
$ws.Range("F7").Value2 = 4.137243980160434
$ws.Range("G7").Value2 = 0.002609482469258721
$ws.Range("J7").Value2 = 0.2339221149679389
$ws.Range("K7").Value2 = 3.910715012255935
$ws.Range("N7").Value2 = 2.154515067465667

$ws.Range("C8").Value2 = 0.04498251807407883
$ws.Range("D8").Value2 = 0.02103136625876267
$ws.Range("E8").Value2 = 0.1065788928573497
$ws.Range("F8").Value2 = 4.196503854875829
$ws.Range("G8").Value2 = 0.002599989064390117
$ws.Range("J8").Value2 = 0.2398669976257679
$ws.Range("K8").Value2 = 4.172062362103134
$ws.Range("N8").Value2 = 2.127761438219949

$ws.Range("C9").Value2 = 0.05793823704993883
$ws.Range("D9").Value2 = 0.02153927250319398
$ws.Range("E9").Value2 = 0.1114154612602256
$ws.Range("F9").Value2 = 4.333726580807763
$ws.Range("G9").Value2 = 0.002583159735728203
$ws.Range("J9").Value2 = 0.2527668912322554
$ws.Range("K9").Value2 = 4.707438833005142
$ws.Range("N9").Value2 = 2.081238579057469

$ws.Range("C10").Value2 = 0.06758704851191055
$ws.Range("D10").Value2 = 0.02193116146774798
$ws.Range("E10").Value2 = 0.115267983985877
$ws.Range("F10").Value2 = 4.447689305572254
$ws.Range("G10").Value2 = 0.002571870167870588
$ws.Range("J10").Value2 = 0.2630225563364661
$ws.Range("K10").Value2 = 5.11528651814757
$ws.Range("N10").Value2 = 2.050728949326327

$ws.Range("C11").Value2 = 0.07200808722718932
$ws.Range("D11").Value2 = 0.02211336164740274
$ws.Range("E11").Value2 = 0.1170866257945207
$ws.Range("F11").Value2 = 4.502463225691315
$ws.Range("G11").Value2 = 0.002566964528057956
$ws.Range("J11").Value2 = 0.2678618121195058
$ws.Range("K11").Value2 = 5.304104344808877
$ws.Range("N11").Value2 = 2.037657195882829

$ws.Range("C12").Value2 = 0.07368703621209249
$ws.Range("D12").Value2 = 0.02218290827958214
$ws.Range("E12").Value2 = 0.1177848858006527
$ws.Range("F12").Value2 = 4.523632220667565
$ws.Range("G12").Value2 = 0.002565139726948053
$ws.Range("J12").Value2 = 0.2697196892441127
$ws.Range("K12").Value2 = 5.376086979131799
$ws.Range("N12").Value2 = 2.032824296587705

$ws.Range("C13").Value2 = 0.07332522738269631
$ws.Range("D13").Value2 = 0.02216790585463357
$ws.Range("E13").Value2 = 0.1176340757536209
$ws.Range("F13").Value2 = 4.519054002004054
$ws.Range("G13").Value2 = 0.002565531272582319
$ws.Range("J13").Value2 = 0.2693184288010855
$ws.Range("K13").Value2 = 5.360562674995094
$ws.Range("N13").Value2 = 2.033859923974077

$ws.Range("C14").Value2 = 0.07214611787100011
$ws.Range("D14").Value2 = 0.02211907231579957
$ws.Range("E14").Value2 = 0.1171438796417164
$ws.Range("F14").Value2 = 4.504196214675801
$ws.Range("G14").Value2 = 0.002566813743407926
$ws.Range("J14").Value2 = 0.26801415058668
$ws.Range("K14").Value2 = 5.310016699124674
$ws.Range("N14").Value2 = 2.037257237828172

$ws.Range("C15").Value2 = 0.07142451101816505
$ws.Range("D15").Value2 = 0.02208923177637345
$ws.Range("E15").Value2 = 0.1168448701677214
$ws.Range("F15").Value2 = 4.495151210089972
$ws.Range("G15").Value2 = 0.002567603566311934
$ws.Range("J15").Value2 = 0.267218555124856
$ws.Range("K15").Value2 = 5.279118793426903
$ws.Range("N15").Value2 = 2.039353469593379

$ws.Range("C16").Value2 = 0.06729878181525351
$ws.Range("D16").Value2 = 0.02191933213745756
$ws.Range("E16").Value2 = 0.1151504679705297
$ws.Range("F16").Value2 = 4.444169180930459
$ws.Range("G16").Value2 = 0.002572195372782637
$ws.Range("J16").Value2 = 0.2627098285919089
$ws.Range("K16").Value2 = 5.103013706041111
$ws.Range("N16").Value2 = 2.051599548674432

$ws.Range("C17").Value2 = 0.06477607546551667
$ws.Range("D17").Value2 = 0.02181610053761318
$ws.Range("E17").Value2 = 0.1141279907270665
$ws.Range("F17").Value2 = 4.413648270422328
$ws.Range("G17").Value2 = 0.002575071059091695
$ws.Range("J17").Value2 = 0.2599886571152865
$ws.Range("K17").Value2 = 4.995826798541941
$ws.Range("N17").Value2 = 2.059319585669954

$ws.Range("C18").Value2 = 0.06332805538801267
$ws.Range("D18").Value2 = 0.02175709494429157
$ws.Range("E18").Value2 = 0.1135461068872061
$ws.Range("F18").Value2 = 4.396369045984585
$ws.Range("G18").Value2 = 0.002576746743536809
$ws.Range("J18").Value2 = 0.2584398583725829
$ws.Range("K18").Value2 = 4.934484458210875
$ws.Range("N18").Value2 = 2.063835865410638

$ws.Range("C19").Value2 = 0.06283828443164907
$ws.Range("D19").Value2 = 0.02173718071267672
$ws.Range("E19").Value2 = 0.1133501567808075
$ws.Range("F19").Value2 = 4.390565758961827
$ws.Range("G19").Value2 = 0.002577317829456473
$ws.Range("J19").Value2 = 0.2579182588681874
$ws.Range("K19").Value2 = 4.913767753277853
$ws.Range("N19").Value2 = 2.065378008719492

$ws.Range("C20").Value2 = 0.06504431240723818
$ws.Range("D20").Value2 = 0.02182705147201247
$ws.Range("E20").Value2 = 0.1142361911106597
$ws.Range("F20").Value2 = 4.416868709251759
$ws.Range("G20").Value2 = 0.002574762696679531
$ws.Range("J20").Value2 = 0.2602766360984248
$ws.Range("K20").Value2 = 5.00720501121333
$ws.Range("N20").Value2 = 2.058489910004027

$ws.Range("C21").Value2 = 0.07249231883569962
$ws.Range("D21").Value2 = 0.02213340104991701
$ws.Range("E21").Value2 = 0.1172876015309043
$ws.Range("F21").Value2 = 4.508548663906367
$ws.Range("G21").Value2 = 0.002566436160723162
$ws.Range("J21").Value2 = 0.2683965577359544
$ws.Range("K21").Value2 = 5.324850139596549
$ws.Range("N21").Value2 = 2.036256178442414

$ws.Range("C22").Value2 = 0.07738809684758508
$ws.Range("D22").Value2 = 0.02233682613322685
$ws.Range("E22").Value2 = 0.1193377491735674
$ws.Range("F22").Value2 = 4.570960195007927
$ws.Range("G22").Value2 = 0.002561185707640018
$ws.Range("J22").Value2 = 0.2738513581896171
$ws.Range("K22").Value2 = 5.535259995672845
$ws.Range("N22").Value2 = 2.022408123873959

$ws.Range("C23").Value2 = 0.07477248102250655
$ws.Range("D23").Value2 = 0.02222796514445591
$ws.Range("E23").Value2 = 0.1182384093297486
$ws.Range("F23").Value2 = 4.537419916804083
$ws.Range("G23").Value2 = 0.002563970531677096
$ws.Range("J23").Value2 = 0.2709263723906332
$ws.Range("K23").Value2 = 5.422700134433057
$ws.Range("N23").Value2 = 2.029736259501107

$ws.Range("C24").Value2 = 0.06492303528537491
$ws.Range("D24").Value2 = 0.02182209948358604
$ws.Range("E24").Value2 = 0.1141872551838148
$ws.Range("F24").Value2 = 4.415411915607251
$ws.Range("G24").Value2 = 0.002574902037339933
$ws.Range("J24").Value2 = 0.2601463921145211
$ws.Range("K24").Value2 = 5.002060047179441
$ws.Range("N24").Value2 = 2.058864763721886

$ws.Range("C25").Value2 = 0.05441156793030189
$ws.Range("D25").Value2 = 0.02139852877850501
$ws.Range("E25").Value2 = 0.1100548630893243
$ws.Range("F25").Value2 = 4.294320331998136
$ws.Range("G25").Value2 = 0.002587522690463137
$ws.Range("J25").Value2 = 0.249142025477525
$ws.Range("K25").Value2 = 4.560099141298053
$ws.Range("N25").Value2 = 2.093183234409196
